# Insert a new row at position 546, shifting existing rows 546:588 down to 547:589,
# then populate the new row 546 with the new daily price entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 546 (rows below shift down by one).
$ws.Rows.Item(546).Insert()

# Populate the newly inserted row 546 with the new record.
$ws.Range("A546").Value = 10
$ws.Range("B546").Value = "Vega Modelo de Temuco"
$ws.Range("C546").Value = "La Araucanía"
$ws.Range("D546").Value = 45223
$ws.Range("E546").Value = 9
$ws.Range("F546").Value = 100112009
$ws.Range("G546").Value = "Acelga"
$ws.Range("H546").Value = "Sin especificar"
$ws.Range("I546").Value = "Primera"
$ws.Range("J546").Value = 40
$ws.Range("K546").Value = 8000
$ws.Range("L546").Value = 8000
$ws.Range("M546").Value = 8000
$ws.Range("N546").Value = "$/docena de atados (12 kilos)"
$ws.Range("O546").Value = "Región de La Araucanía"
$ws.Range("P546").Value = 667
$ws.Range("Q546").Value = 12
$ws.Range("R546").Value = "Hortaliza"
